# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Update the daily conversion note on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.04 = 27713.22 pesos`n✅ 27713.22 pesos = 7.0 = 966.29 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate cells on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 141.99
$wsTasas.Range("O10").Value = 3935
$wsTasas.Range("N12").Value = 3959.99
$wsTasas.Range("O12").Value = 138.075
